$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 466.33334
$ws.Range("I32").Value = 450
$ws.Range("J32").Value = 499
$ws.Range("K32").Value = 450
$ws.Range("L32").Value = 499
$ws.Range("M32").Value = -124
$ws.Range("N32").Value = -1151
$ws.Range("H64").Value = 4584.421
$ws.Range("J64").Value = 3500
$ws.Range("L64").Value = 3500
$ws.Range("N64").Value = -3996
$ws.Range("H67").Value = 4584.421
$ws.Range("J67").Value = 3500
$ws.Range("L67").Value = 3500
$ws.Range("N67").Value = -5216
$ws.Range("H74").Value = 4725
$ws.Range("I74").Value = 5000
$ws.Range("J74").Value = 3900
$ws.Range("K74").Value = 5000
$ws.Range("L74").Value = 3900
$ws.Range("M74").Value = -4064
$ws.Range("N74").Value = -5772
$ws.Range("H76").Value = 4508743
$ws.Range("I76").Value = 8775283
$ws.Range("K76").Value = 8775283
$ws.Range("M76").Value = -8774968
$ws.Range("H77").Value = 4725
$ws.Range("I77").Value = 5000
$ws.Range("J77").Value = 3900
$ws.Range("K77").Value = 25000
$ws.Range("L77").Value = 19500
$ws.Range("M77").Value = -20320
$ws.Range("N77").Value = -28860
$ws.Range("H79").Value = 4508743
$ws.Range("I79").Value = 8775283
$ws.Range("K79").Value = 8775283
$ws.Range("M79").Value = -8774191
$ws.Range("H81").Value = 36000
$ws.Range("J81").Value = 36000
$ws.Range("L81").Value = 36000
$ws.Range("N81").Value = -37996
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H84").Value = 36000
$ws.Range("J84").Value = 36000
$ws.Range("L84").Value = 108000
$ws.Range("N84").Value = -117984
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H86").Value = 4278.4346
$ws.Range("I86").Value = 4421.8237
$ws.Range("K86").Value = 4421.8237
$ws.Range("M86").Value = -3298.8237
$ws.Range("H87").Value = 29899.5
$ws.Range("J87").Value = 29899.5
$ws.Range("L87").Value = 29899.5
$ws.Range("N87").Value = -32395.5
$ws.Range("H89").Value = 4278.4346
$ws.Range("I89").Value = 4421.8237
$ws.Range("K89").Value = 22109.1185
$ws.Range("M89").Value = -16493.1185
$ws.Range("H90").Value = 29899.5
$ws.Range("J90").Value = 29899.5
$ws.Range("L90").Value = 89698.5
$ws.Range("N90").Value = -102178.5
$ws.Range("H92").Value = 81699840
$ws.Range("I92").Value = 3704266.8
$ws.Range("K92").Value = 3704266.8
$ws.Range("M92").Value = -3703018.8
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H97").Value = 933.3333
$ws.Range("I97").Value = 600
$ws.Range("J97").Value = 1028.5714
$ws.Range("K97").Value = 1800
$ws.Range("L97").Value = 3085.7142
$ws.Range("M97").Value = -1304
$ws.Range("N97").Value = -4077.7142
$ws.Range("H99").Value = 1566
$ws.Range("I99").Value = 632
$ws.Range("J99").Value = 2500
$ws.Range("K99").Value = 1896
$ws.Range("L99").Value = 7500
$ws.Range("M99").Value = -398
$ws.Range("N99").Value = -10496
$ws.Range("H141").Value = 2071.9473
$ws.Range("I141").Value = 1197.4286
$ws.Range("J141").Value = 4520.6
$ws.Range("K141").Value = 3592.2858
$ws.Range("L141").Value = 13561.8
$ws.Range("M141").Value = 1587.7142
$ws.Range("N141").Value = -23921.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 4632517
$ws.Range("I102").Value = 6174672.5
$ws.Range("J102").Value = 6050
$ws.Range("K102").Value = 6174672.5
$ws.Range("L102").Value = 6050
$ws.Range("M102").Value = -6173050.5
$ws.Range("N102").Value = -9294
$ws.Range("H110").Value = 14683.826
$ws.Range("I110").Value = 17359.37
$ws.Range("J110").Value = 1975
$ws.Range("K110").Value = 17359.37
$ws.Range("L110").Value = 1975
$ws.Range("M110").Value = -15314.37
$ws.Range("N110").Value = -6065
$ws.Range("H135").Value = 53692
$ws.Range("J135").Value = 53692
$ws.Range("L135").Value = 53692
$ws.Range("N135").Value = -63832

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2687
$ws.Range("I105").Value = 2675
$ws.Range("K105").Value = 2675
$ws.Range("M105").Value = -928
$ws.Range("H107").Value = 216591.64
$ws.Range("I107").Value = 302607
$ws.Range("J107").Value = 1553.25
$ws.Range("K107").Value = 302607
$ws.Range("L107").Value = 1553.25
$ws.Range("M107").Value = -300687
$ws.Range("N107").Value = -5393.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 7609.8096
$ws.Range("I62").Value = 8271.429
$ws.Range("J62").Value = 6286.5713
$ws.Range("K62").Value = 8271.429
$ws.Range("L62").Value = 6286.5713
$ws.Range("M62").Value = -7647.429
$ws.Range("N62").Value = -7534.5713
$ws.Range("H65").Value = 7609.8096
$ws.Range("I65").Value = 8271.429
$ws.Range("J65").Value = 6286.5713
$ws.Range("K65").Value = 41357.145
$ws.Range("L65").Value = 31432.8565
$ws.Range("M65").Value = -38237.145
$ws.Range("N65").Value = -37672.85649999999
$ws.Range("H105").Value = 1887.4193
$ws.Range("I105").Value = 1977.3077
$ws.Range("J105").Value = 1420
$ws.Range("K105").Value = 1977.3077
$ws.Range("L105").Value = 1420
$ws.Range("M105").Value = -230.3077000000001
$ws.Range("N105").Value = -4914
$ws.Range("H132").Value = 5407735.5
$ws.Range("I132").Value = 9092475
$ws.Range("J132").Value = 3450
$ws.Range("K132").Value = 27277425
$ws.Range("L132").Value = 10350
$ws.Range("M132").Value = -27274895
$ws.Range("N132").Value = -15410

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 15879541
$ws.Range("I117").Value = 25389.5
$ws.Range("J117").Value = 19609930
$ws.Range("K117").Value = 76168.5
$ws.Range("L117").Value = 58829790
$ws.Range("M117").Value = -72726.5
$ws.Range("N117").Value = -58836674
$ws.Range("H131").Value = 2858178.5
$ws.Range("I131").Value = 5263727
$ws.Range("J131").Value = 1590
$ws.Range("K131").Value = 15791181
$ws.Range("L131").Value = 4770
$ws.Range("M131").Value = -15786141
$ws.Range("N131").Value = -14850

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 500000
$ws.Range("I26").Value = 500000
$ws.Range("K26").Value = 500000
$ws.Range("M26").Value = -499720
$ws.Range("H50").Value = 500000
$ws.Range("I50").Value = 500000
$ws.Range("K50").Value = 500000
$ws.Range("M50").Value = -499502
$ws.Range("H80").Value = 10637.917
$ws.Range("J80").Value = 2500
$ws.Range("L80").Value = 2500
$ws.Range("N80").Value = -4496
$ws.Range("H83").Value = 10637.917
$ws.Range("J83").Value = 2500
$ws.Range("L83").Value = 12500
$ws.Range("N83").Value = -22484
$ws.Range("H102").Value = 2715.3225
$ws.Range("I102").Value = 2564.1304
$ws.Range("J102").Value = 3150
$ws.Range("K102").Value = 2564.1304
$ws.Range("L102").Value = 3150
$ws.Range("M102").Value = -942.1304
$ws.Range("N102").Value = -6394

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 83336056
$ws.Range("I40").Value = 111113690
$ws.Range("K40").Value = 111113690
$ws.Range("M40").Value = -111113554
$ws.Range("H125").Value = 45000
$ws.Range("J125").Value = 45000
$ws.Range("L125").Value = 45000
$ws.Range("N125").Value = -54840
$ws.Range("H135").Value = 133251.14
$ws.Range("J135").Value = 133251.14
$ws.Range("L135").Value = 133251.14
$ws.Range("N135").Value = -143391.14
$ws.Range("H138").Value = 68865.42999999999
$ws.Range("J138").Value = 68865.42999999999
$ws.Range("L138").Value = 68865.42999999999
$ws.Range("N138").Value = -79145.42999999999
$ws.Range("H140").Value = 67892.08
$ws.Range("J140").Value = 67892.08
$ws.Range("L140").Value = 67892.08
$ws.Range("N140").Value = -78252.08

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 71429350
$ws.Range("I107").Value = 166667340
$ws.Range("J107").Value = 862.5
$ws.Range("K107").Value = 500002020
$ws.Range("L107").Value = 2587.5
$ws.Range("M107").Value = -500000100
$ws.Range("N107").Value = -6427.5
$ws.Range("H126").Value = 1225.4
$ws.Range("I126").Value = 1090.4615
$ws.Range("J126").Value = 2102.5
$ws.Range("K126").Value = 3271.3845
$ws.Range("L126").Value = 6307.5
$ws.Range("M126").Value = -801.3844999999997
$ws.Range("N126").Value = -11247.5

